# Apply the edit described by the commit:
#   "Updated readme correct worksheet name"
#
# Logical changes:
#   1. Rename sheet 2 "User specific label" -> "Multilingual icon"
#   2. Fix the instruction text in that sheet (B6): the referenced cell
#      moved from A14 to A16, so update the label accordingly.
#   3. Update view-state: the "Multilingual icon" sheet becomes the active
#      tab/selection, with the other sheets' selections moved along as
#      recorded when the workbook was last saved.

$wb = $excel.ActiveWorkbook

$wsChoice = $wb.Worksheets.Item("Choice field")
$wsLabel  = $wb.Worksheets.Item("User specific label")
$wsSql    = $wb.Worksheets.Item("SQL Command templates")

# 1. Rename the worksheet
$wsLabel.Name = "Multilingual icon"

# 2. Correct the instruction text (cell reference A14 -> A16)
$wsLabel.Range("B6").Value = "Copy case statement below the table in A16"

# 3. Restore/update selections on each sheet
$wsChoice.Range("A19").Select() | Out-Null
$wsSql.Range("B2").Select() | Out-Null

# Make the renamed sheet the active tab, with its own selection
$wsLabel.Activate()
$wsLabel.Range("B20").Select() | Out-Null
